$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.214426040649414
$ws.Range("B1").Value = 2.445568323135376
$ws.Range("C1").Value = 4.807501316070557
$ws.Range("D1").Value = 2.527350664138794
$ws.Range("E1").Value = 1.079631447792053
